# Logged Week 15 and simulated Week 16
# Update the "R" (Road) row target-depth totals on both the OFF and DEF
# sheets to reflect the newly logged/simulated week's cumulative numbers.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 ("R") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 254
$wsOff.Range("C3").Value = 182
$wsOff.Range("D3").Value = 133
$wsOff.Range("E3").Value = 54
$wsOff.Range("F3").Value = 4

# --- DEF sheet: row 3 ("R") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 419
$wsDef.Range("C3").Value = 294
$wsDef.Range("D3").Value = 116
$wsDef.Range("E3").Value = 43
$wsDef.Range("F3").Value = 10
